# Append the newest COVID-19 daily data rows (25-31 May 2020) to the
# "Covid-19 podatki" sheet / "Tabela1" table, matching the upstream
# GitHub-bot data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19 podatki")
$lo = $ws.ListObjects.Item(1)

# Date (serial), Tested (all), Tested (daily), Positive (all), Positive (daily),
# All hospitalized, All in intensive care, Discharged, Deaths (all), Deaths (daily)
$data = @(
    @(43976, 75770, 754, 1469, 0, 9, 2, 6, 108, 1),
    @(43977, 76579, 809, 1471, 2, 8, 2, 2, 108, 0),
    @(43978, 77210, 631, 1473, 2, 7, 2, 1, 108, 0),
    @(43979, 77916, 706, 1473, 0, 7, 2, 0, 108, 0),
    @(43980, 78529, 613, 1473, 0, 7, 2, 0, 108, 0),
    @(43981, 78793, 264, 1473, 0, 6, 2, 1, 108, 0),
    @(43982, 79039, 246, 1473, 0, 5, 1, 0, 109, 1)
)

# Grow the sheet/table one row at a time by duplicating the last existing
# row (copy + insert) so every new row inherits the same cell formatting
# (date format / thousands separator / general) already used by the table,
# then overwrite the copied values with the real data for that day.
$lastRow = $ws.Range("A1048576").End(-4162).Row
foreach ($rowVals in $data) {
    $newRowNum = $lastRow + 1
    $ws.Rows($lastRow).Copy()
    $ws.Rows($newRowNum).Insert(-4121)

    for ($col = 1; $col -le 10; $col++) {
        $ws.Cells.Item($newRowNum, $col).Value = $rowVals[$col - 1]
    }

    $lastRow = $newRowNum
}

# Resize the table (and its autofilter) to cover the newly added rows.
$lo.Resize($ws.Range("A1:J" + $lastRow))

# Reflect the new last row in the sheet's selection, as a user scrolling
# down to review the freshly entered data would leave it.
$ws.Activate()
$ws.Range("A" + $lastRow + ":J" + $lastRow).Select()
